$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new values look numeric must be forced to Text format first,
# so Excel stores them as strings (matching the original inlineStr data)
# instead of auto-converting them to numbers.
$textCells = @(
    "D5", "D6", "D7", "D8", "D9", "D10", "D11", "D12", "D13", "D14", "D15", "D16", "D17", "D18", "D20", "D23", "D24", "D25", "D26", "D27", "D28", "D29", "D30", "D31", "D32", "D33", "D34", "D36", "D37", "D39", "D40", "D41", "D42", "D43", "D44", "D45", "D46", "D47", "D48", "D49", "D50", "D51"
)
foreach ($cellRef in $textCells) {
    $ws.Range($cellRef).NumberFormat = "@"
}

$ws.Range("D2").Value = '27.207.61'
$ws.Range("E2").Value = '  -0.09%  '

$ws.Range("D3").Value = '1.903.57'
$ws.Range("E3").Value = '  +0.38%  '

$ws.Range("E4").Value = '  -0.13%  '

$ws.Range("D5").Value = '306.06'
$ws.Range("E5").Value = '  -0.52%  '

$ws.Range("D6").Value = '0.9999'
$ws.Range("E6").Value = '  -0.20%  '

$ws.Range("D7").Value = '0.5365'
$ws.Range("E7").Value = '  +3.40%  '

$ws.Range("D8").Value = '0.3808'
$ws.Range("E8").Value = '  +1.27%  '

$ws.Range("D9").Value = '0.07283'
$ws.Range("E9").Value = '  -0.04%  '

$ws.Range("D10").Value = '22.16'
$ws.Range("E10").Value = '  +4.56%  '

$ws.Range("D11").Value = '0.9049'
$ws.Range("E11").Value = '  +0.53%  '

$ws.Range("D12").Value = '0.08195'
$ws.Range("E12").Value = '  +0.68%  '

$ws.Range("D13").Value = '95.90'
$ws.Range("E13").Value = '  -0.51%  '

$ws.Range("D14").Value = '5.337'
$ws.Range("E14").Value = '  +1.12%  '

$ws.Range("D15").Value = '1.001'
$ws.Range("E15").Value = '  -0.08%  '

$ws.Range("D16").Value = '14.85'
$ws.Range("E16").Value = '  +2.16%  '

$ws.Range("D17").Value = '0.000008646'
$ws.Range("E17").Value = '  +0.34%  '

$ws.Range("D18").Value = '0.9999'
$ws.Range("E18").Value = '  -0.18%  '

$ws.Range("D19").Value = '27.230.29'
$ws.Range("E19").Value = '  -0.10%  '

$ws.Range("D20").Value = '5.040'
$ws.Range("E20").Value = '  -0.89%  '

$ws.Range("D21").Value = '1.078.62'
$ws.Range("E21").Value = '  -43.14%  '

$ws.Range("E22").Value = '  +0.86%  '

$ws.Range("D23").Value = '6.511'
$ws.Range("E23").Value = '  +1.81%  '

$ws.Range("D24").Value = '149.29'
$ws.Range("E24").Value = '  +1.49%  '

$ws.Range("D25").Value = '2.291'

$ws.Range("D26").Value = '18.35'
$ws.Range("E26").Value = '  +0.78%  '

$ws.Range("D27").Value = '1.748'
$ws.Range("E27").Value = '  +0.26%  '

$ws.Range("D28").Value = '116.65'
$ws.Range("E28").Value = '  +1.23%  '

$ws.Range("D29").Value = '4.812'
$ws.Range("E29").Value = '  -0.37%  '

$ws.Range("D30").Value = '4.727'
$ws.Range("E30").Value = '  -4.57%  '

$ws.Range("D31").Value = '0.09211'
$ws.Range("E31").Value = '  -0.25%  '

$ws.Range("D32").Value = '0.8294'
$ws.Range("E32").Value = '  +4.33%  '

$ws.Range("D33").Value = '0.05074'
$ws.Range("E33").Value = '  +0.74%  '

$ws.Range("D34").Value = '1.217'
$ws.Range("E34").Value = '  +0.07%  '

$ws.Range("E35").Value = '  +1.72%  '

$ws.Range("D36").Value = '3.340'
$ws.Range("E36").Value = '  -3.22%  '

$ws.Range("D37").Value = '2.676'
$ws.Range("E37").Value = '  +3.08%  '

$ws.Range("D39").Value = '0.02003'

$ws.Range("D40").Value = '1.075'
$ws.Range("E40").Value = '  +0.18%  '

$ws.Range("D41").Value = '9.321'
$ws.Range("E41").Value = '  +4.03%  '

$ws.Range("D42").Value = '6.617'
$ws.Range("E42").Value = '  +0.89%  '

$ws.Range("D43").Value = '117.11'
$ws.Range("E43").Value = '  +1.55%  '

$ws.Range("B44").Value = 'Decentraland'
$ws.Range("C44").Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range("D44").Value = '0.5036'
$ws.Range("E44").Value = '  +3.68%  '

$ws.Range("B45").Value = 'Algorand'
$ws.Range("C45").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D45").Value = '0.1521'
$ws.Range("E45").Value = '  +0.45%  '

$ws.Range("D46").Value = '0.9991'
$ws.Range("E46").Value = '  -0.31%  '

$ws.Range("D47").Value = '10.14'
$ws.Range("E47").Value = '  +1.25%  '

$ws.Range("D48").Value = '1.638'
$ws.Range("E48").Value = '  +1.06%  '

$ws.Range("D49").Value = '38.35'
$ws.Range("E49").Value = '  +0.38%  '

$ws.Range("D50").Value = '0.06166'
$ws.Range("E50").Value = '  +3.78%  '

$ws.Range("D51").Value = '63.38'
$ws.Range("E51").Value = '  -0.02%  '

# Restore default (Normal) style on the cells we forced to Text format,
# so no stray style index is left behind on them.
foreach ($cellRef in $textCells) {
    $ws.Range($cellRef).Style = "Normal"
}
